$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C11"  = -13.0427
    "A12"  = -21.39289999999999
    "C23"  = -12.08909999999999
    "A27"  = -21.8035
    "C28"  = -13.6107
    "A32"  = -21.515
    "C32"  = -12.8049
    "C34"  = -12.16160000000001
    "A36"  = -19.9274
    "A38"  = -19.50389999999998
    "C42"  = -12.5894
    "A46"  = -21.89120000000002
    "C49"  = -13.6343
    "A54"  = -21.75319999999998
    "C54"  = -12.8655
    "A55"  = -22.2644
    "A56"  = -22.05260000000003
    "A67"  = -21.43979999999997
    "A69"  = -21.73239999999997
    "A72"  = -21.8478
    "C78"  = -13.2212
    "C80"  = -13.09890000000001
    "A83"  = -21.54699999999998
    "A86"  = -21.90640000000001
    "A91"  = -20.54249999999997
    "A93"  = -21.36330000000001
    "C97"  = -11.87510000000001
    "A99"  = -21.7894
    "C99"  = -12.5886
    "C101" = -12.813
    "A104" = -21.50629999999999
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
